$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each Price/Volume cell stores text (e.g. "300.35", "-0.29%"), not a
# number, so force Text number format before assigning the literal string
# -- otherwise Excel auto-converts it to a numeric/percentage value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.78%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.951"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.59%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07677"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.973"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-15.49%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.829"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.15%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.800"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.97%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9200"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.11%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1753"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.31%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.53%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08615"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.03%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03180"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.12%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.13%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001520"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.06%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005784"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.51%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.26%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.151"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.27%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.35%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1327"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.17%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.277"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.09%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1994"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.38%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04522"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.02%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001224"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004410"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.37%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001253"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.27%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01695"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04681"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007516"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.76%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1350"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.64%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002336"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.68%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.17%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006255"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.51%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8234"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-28.52%"
